$d = $word.ActiveDocument

# --- 1. Title ---
$d.Content.Find.Execute("Nanobots: The Future of Medicine", $true, $false, $false, $false, $false, $true, 1, $false, "The Heartbeat of Healing: A Journey into Medicine", 2) | Out-Null

# --- 2. Author name paragraph: "Dr. Katherine Abernathy" -> "Amelia Harrison" ---
$pAuthor = $d.Paragraphs.Item(2)
$rAuthor = $d.Range($pAuthor.Range.Start, $pAuthor.Range.End - 1)
$rAuthor.Text = "Amelia Harrison"

# --- 3. Email paragraph: "katherine.abernathy@healthcare.edu" -> "ameliaharr1021@gmail.com" ---
$pEmail = $d.Paragraphs.Item(3)
$rEmail = $d.Range($pEmail.Range.Start, $pEmail.Range.End - 1)
$rEmail.Text = "ameliaharr1021@gmail.com"

# --- 4. Body paragraph: full rewrite (keeps internal line breaks via vertical-tab) ---
$pBody = $d.Paragraphs.Item(5)
$rBody = $d.Range($pBody.Range.Start, $pBody.Range.End - 1)
$rBody.Text = "As we embark on this chapter of learning and exploration, we will delve into the realm of medicine, an extraordinary field dedicated to preserving and restoring human health. From the intricate harmony of our physiology to the boundless possibilities of modern therapies, we will discover the beauty and complexities of medicine.`v`vThe human body, a marvel of biological engineering, serves as the canvas upon which medicine works its magic. We will journey through the interconnected systems that govern our existence, from the microscopic world of cells to the intricate network of organs. We will unravel the secrets of human physiology, discovering the delicate balance that maintains our health and the ways in which medicine can intervene when this balance is disrupted.`v`vAt the heart of medicine lies the patient, an individual with unique stories, hopes, and fears. As we step into the world of healing, we will explore the art of patient care, emphasizing empathy, compassion, and respect for human dignity. We will learn how medicine goes beyond treating symptoms; it also involves nurturing the human spirit, fostering a bond of trust that empowers patients to actively participate in their healing journey.`v`vIntroduction Continued:`v`vThe practice of medicine draws upon a vast reservoir of knowledge, including scientific research, clinical experience, and cultural traditions. We will examine how evidence-based medicine guides medical practice, ensuring that treatments are safe, effective, and tailored to individual patient needs. We will also explore alternative and complementary therapies, recognizing the diversity of approaches to healing.`v`vAs medicine evolves, it faces numerous challenges, from emerging diseases and antibiotic resistance to the complexities of healthcare systems and the rising cost of medical care. We will delve into these issues, seeking a deeper understanding of the intricate factors that shape modern medicine. Through critical thinking and open-minded discussions, we will explore potential solutions and envision a future where medical advancements benefit all members of society.`v`vIntroduction Concluded:`v`vMedicine is a noble profession, one that intertwines science, art, and unwavering dedication to patient care. As we traverse the landscape of healing, may we cultivate a deep appreciation for the complex tapestry of human life, the indomitable power of the human spirit, and the profound responsibility we carry as future guardians of health."

# --- 5. Summary paragraph text rewrite ---
$pSummary = $d.Paragraphs.Item(7)
$rSummary = $d.Range($pSummary.Range.Start, $pSummary.Range.End - 1)
$rSummary.Text = "Our exploration of medicine has unveiled the profound impact it wields on human lives, delving into the remarkable complexity of the human body, the challenges of modern medicine, and the essential role of empathy and respect in patient care. We have unraveled the delicate interplay between science, tradition, and cultural factors that shape medical practices. Ultimately, we have gained an appreciation for the profound responsibility that rests upon those dedicated to preserving and restoring human health."

# --- 6. New empty paragraph at the very end of the document body ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
